$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "Unnamed: 0" column (C); D:H shift left to C:G, preserving styles
$ws.Range("C1").EntireColumn.Delete()

# Remove old transaction rows 3-6, keeping header (row1) + single data row (row2)
$ws.Range("A3:A6").EntireRow.Delete()

# Fix header text in B1
$ws.Range("B1").Value = "Unnamed: 0"

# Update the remaining data row (row 2) with the new transaction values
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = "Direct Deposit"
$ws.Range("E2").Value = 5000
$ws.Range("F2").Value = 0

# G2 holds a date-like string that must stay plain text, not be parsed as a date
$g2 = $ws.Range("G2")
$g2.NumberFormat = "@"
$g2.Value = "1/13/2025"
$g2.NumberFormat = "General"
$g2.Style = "Normal"
